$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document ("... each list
# has a lock.") -- it currently carries the _GoBack bookmark. We are
# going to append four new paragraphs after it and move the _GoBack
# bookmark onto the very last one of them, exactly like Word does when
# you keep typing at the end of a document.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$anchorEnd = $lastPara.Range.End

# Remove the existing _GoBack bookmark -- it will be re-created at the
# new end of the document once the new paragraphs are in place.
if ($d.Bookmarks.Exists("_GoBack")) {
    [void]$d.Bookmarks.Item("_GoBack").Delete()
}

# Insert right before the paragraph mark of the last paragraph (i.e. one
# character before the absolute end of the story) so the new content is
# appended immediately after the existing text instead of replacing it.
$insertionPoint = $d.Range($anchorEnd - 1, $anchorEnd - 1)

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>22/03/14</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>I decided to put the resize function on hold and divert my attention to modifying the locked list so that each list had its own lock, so when a thread entered a critical section it need only lock the list it was in, and not the entire table, as with hash_locked.cpp.</w:t></w:r><w:r><w:t xml:space="preserve"> This proved relatively easy</w:t></w:r><w:r><w:t xml:space="preserve"> as, instead of locking down the whole table, I instead just locked by list. I now must compare both hash_locked and hash_locked_per_bucket to see if there is a performance difference.</w:t></w:r><w:r><w:t xml:space="preserve"> Initially there does not appear to be a major difference when using plain mutex locks, though per bucket seems to have a slight advantage at higher thread counts</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>I have attempted to add a resize functionality to the lockless hash table but as of yet I have been unsuccessful, I feel that this may be a step too far and that I may need to leave it unimplemented, or failing that, implement the locked resize function that I have already.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Started work on my presentation slides.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

[void]$insertionPoint.InsertXML($newParagraphsXml)
